# Update the "Program_choosing" worksheet with the new list of programs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# New ordered list of program names for column A (rows 2-15).
$programs = @(
    "TUM_Informatics",
    "RWTH_Aachen_Data Science_CS_BG",
    "RWTH_Aachen_Data Science_Math_BG",
    "RWTH_Aachen_Software_System_Engineering",
    "RWTH_Aachen_Media_Informatics",
    "Freie Uni Berlin - Data Science",
    "TU Berlin Computer Science",
    "TUM Data Engineering and Analytics",
    "TU Delft Computer Science",
    "RWTH_Aachen_DDS",
    "RWTH_Aachen_TIME",
    "Uni_Goettingen_Applied_CS",
    "TUM_Math_Data_Science_MathBackground",
    "TUM_Math_Data_Science_CSBackground"
)

$row = 2
foreach ($program in $programs) {
    $ws.Cells.Item($row, 1).Value = $program
    $ws.Cells.Item($row, 2).Value = "Yes"
    $row = $row + 1
}

$lastRow = $row - 1

# Refresh the data validation list so it covers the new row range (B1:B15).
$validationRange = $ws.Range("B1:B$lastRow")
$validationRange.Validation.Delete()
$validationRange.Validation.Add(3, 1, 1, '"Yes,No"')
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

# Update the selected cell to match the diff (A5 selected).
$ws.Range("A5").Select()
